# Update Name of Algo
# Apply updated imputed values produced by the KNN algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3"   = 5.623
    "C3"   = -12.313
    "D8"   = -8.401
    "D11"  = -7.457000000000001
    "A12"  = -21.629
    "B14"  = 6.228000000000001
    "D14"  = -7.678000000000002
    "D15"  = -8.175000000000001
    "D17"  = -8.177999999999999
    "C20"  = -12.149
    "C25"  = -12.738
    "B26"  = 6.593999999999999
    "D26"  = -8.034000000000001
    "A27"  = -21.217
    "C30"  = -12.717
    "B31"  = 6.276999999999999
    "A32"  = -21.471
    "B35"  = 7.342000000000001
    "A36"  = -20.761
    "D36"  = -8.375
    "B37"  = 7.641
    "A38"  = -20.363
    "C44"  = -12.646
    "B45"  = 5.673
    "A46"  = -21.515
    "C47"  = -12.392
    "B52"  = 5.4
    "A54"  = -21.638
    "A55"  = -21.825
    "A56"  = -21.832
    "B57"  = 6.090000000000001
    "C58"  = -12.802
    "D64"  = -7.616
    "A67"  = -21.588
    "A69"  = -21.721
    "A72"  = -21.567
    "C78"  = -12.536
    "D79"  = -7.813
    "B81"  = 6.430000000000001
    "A83"  = -20.943
    "B83"  = 7.357000000000001
    "C84"  = -13.112
    "A86"  = -22.325
    "C89"  = -11.26
    "D89"  = -6.804
    "A91"  = -21.747
    "C91"  = -11.21
    "C92"  = -11.539
    "A93"  = -21.49
    "C96"  = -13.112
    "A99"  = -20.683
    "B100" = 6.218999999999999
    "B102" = 6.667
    "C102" = -12.798
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
